# Applies the target edit to the workbook:
#  - Extends Sheet1 with 22 new rows (73-94), containing additional
#    simulation results for the "1000P periods" run configuration.
#  - The new dimension becomes A1:W94 (from A1:W72).
#  - xl/sharedStrings.xml grows accordingly as new unique text values
#    (including "1000P periods" and a batch of numeric-array strings)
#    are introduced by these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Sheet1')

$startRow = 73

$newRows = @(
    @(8, 1.5, 1, 5, 0, 10, 0, 0, 0, '1000P periods', 0.01, '[-6.046939064287374e-11,2.970806695705539e-12]', '[[0.0,0.0,975.641655082653],[0.0,0.0,975.6416551017463]]', '[-3614.95414535744,-2883.718844549143]', '[-1337.228393030528,-829.3502517782138]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.12307000160217285, '[0.0004563137537233886,0.0008053631485530542]', 0, 73, 1, 924),
    @(8, 1.5, 1, 5, 0, 11, 0, 0, 0, '1000P periods', 0.01, '[-3.007353937740412e-11,1.3955310207952258e-9]', '[[0.0,0.0,1000.2004065221081],[0.0,0.0,1000.2004069523357]]', '[-3565.693567912392,-2815.9767215302654]', '[-1202.920629494612,-722.5493307063912]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 1.2333011627197266, '[0.0004564579821494023,0.0008063666218747724]', 0, 74, 1, 65728),
    @(8, 5, 1, 25, 0, 100, 0, 0, 0, '1000P periods', 0.01, '[0.0,1.2456245400627403e-7]', '[[0.0,0.0,4896.126236248317],[0.0,0.0,0.0]]', '[-2293.985269926886,-2260.6389730146825]', '[-169.95291397371298,-65.60678614472278]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 79.03585505485535, '[0.00626711358588136,0.006471667054336739]', 0, 75, 1, 4854461),
    @(8, 5, 1, 25, 0, 100, 0, 18, 0, '1000P periods', 0.01, '[-1.2669332101872325e-11,9.827117504271895e-11]', '[[0.0,0.0,4536.001800930624],[0.0,0.0,0.0]]', '[-2515.245289129506,-2260.7145887270026]', '[-176.41908955386344,123.6320244138613]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.11619305610656738, '[0.004259098964354781,0.007403528731528524]', 0, 76, 1, 4752),
    @(8, 5, 1, 25, 0, 100, 0, 36, 0, '1000P periods', 0.01, '[0.0,2.5104293858477443e-11]', '[[0.0,0.0,3593.183258829851],[0.0,0.0,0.0]]', '[-2291.7976697123177,-2265.004858314624]', '[-159.37230699154404,-81.68773353013063]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.02109217643737793, '[0.006352710121004197,0.006432902849153114]', 0, 77, 1, 904),
    @(8, 5, 1, 25, 0, 100, 0, 54, 0, '1000P periods', 0.01, '[0.0,2.218291291511624e-11]', '[[0.0,0.0,2427.795450191914],[0.0,0.0,0.0]]', '[-2289.653212833143,-2266.9510030633082]', '[-148.2255842134629,-79.02619095378725]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.019444942474365234, '[0.006366689520028675,0.006445856707297158]', 0, 78, 1, 808),
    @(8, 5, 1, 25, 0, 100, 0, 72, 0, '1000P periods', 0.01, '[0.0,2.466204295477768e-11]', '[[0.0,0.0,1484.9769079579821],[0.0,0.0,0.0]]', '[-2287.6028431838313,-2262.423348552621]', '[-150.6043089222757,-73.46478623816222]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.028664112091064453, '[0.006367392989904797,0.0064562850847207755]', 0, 79, 1, 902),
    @(8, 5, 1, 25, 0, 100, 0, 90, 0, '1000P periods', 0.01, '[0.0,1.2377466136566902e-7]', '[[0.0,0.0,1124.852068124762],[0.0,0.0,0.0]]', '[-2285.197192732333,-2253.91836480807]', '[-169.52755173827566,-65.42887438694277]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 108.81562614440918, '[0.00634722009269794,0.006470293908743449]', 0, 80, 1, 4826105),
    @(8, 5, 1, 25, 0, 100, 0, 108, 0, '1000P periods', 0.01, '[-6.0498831134854026e-12,3.547338094843459e-11]', '[[0.0,0.0,1484.9769079164403],[0.0,0.0,0.0]]', '[-2321.5754937581646,-2217.0744717227403]', '[-201.97002151694699,6.878024739128932]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.0507969856262207, '[0.0043391092500468555,0.006690378662341604]', 0, 81, 1, 1886),
    @(8, 5, 1, 25, 0, 100, 0, 126, 0, '1000P periods', 0.01, '[-9.150420288664775e-10,9.477443016957469e-11]', '[[0.0,0.0,2427.7954501909057],[0.0,0.0,0.0]]', '[-3420.1231167807955,-1823.9820839116078]', '[-1248.070538397003,5722.8847488173415]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.04879403114318848, '[0.0001737438445264283,0.010567012128438384]', 0, 82, 1, 1911),
    @(8, 5, 1, 25, 0, 100, 0, 144, 0, '1000P periods', 0.01, '[-1.186851579098434e-5,6.02514504985459e-12]', '[[0.0,0.0,3593.1832565132036],[0.0,0.0,3593.18325926986]]', '[-3901.034616056694,3.566423595983679e6]', '[-2738.706803962937,6.271900045101356e7]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.40076303482055664, '[2.771850425589169e-10,0.028948795598466654]', 0, 83, 1, 16704),
    @(8, 5, 1, 25, 0, 100, 0, 162, 0, '1000P periods', 0.01, '[-7.831839217506987e-11,1.0010993407217755e-10]', '[[0.0,0.0,4536.001801017598],[0.0,0.0,4536.001801139413]]', '[-2328.1472858283423,-2138.680316676339]', '[-189.86396520111043,-66.40182835263852]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.06570792198181152, '[0.003790278375165637,0.010040638685672277]', 0, 84, 1, 2378),
    @(8, 5, 1, 25, 0, 100, 0, 180, 0, '1000P periods', 0.01, '[0.0,1.2448868410042784e-7]', '[[0.0,0.0,4896.126236312554],[0.0,0.0,0.0]]', '[-2292.988654749074,-2260.1833513537367]', '[-170.10631252411568,-65.91486758054157]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 110.05816006660461, '[0.0062884227062314745,0.00645679802781983]', 0, 85, 1, 4852042),
    @(8, 5, 1, 25, 0, 100, 0, 198, 0, '1000P periods', 0.01, '[-4.109655657185108e-11,9.70258706813093e-11]', '[[0.0,0.0,4536.001800926036],[0.0,0.0,0.0]]', '[-2541.344407412897,-2224.2828633117965]', '[-199.7164296243186,268.16226216064877]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.08002805709838867, '[0.003958283722565165,0.007558335328314811]', 0, 86, 1, 4717),
    @(8, 5, 1, 25, 0, 100, 0, 216, 0, '1000P periods', 0.01, '[0.0,2.3899284405576858e-11]', '[[0.0,0.0,3593.1832588324605],[0.0,0.0,0.0]]', '[-2289.46216649518,-2269.3968472689203]', '[-153.28654634529408,-77.84794638968998]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.015024900436401367, '[0.0063565851873717,0.006427074734661797]', 0, 87, 1, 890),
    @(8, 5, 1, 25, 0, 100, 0, 234, 0, '1000P periods', 0.01, '[-5.354094341634495e-16,2.1845954202548454e-11]', '[[0.0,0.0,2427.7954501923646],[0.0,0.0,0.0]]', '[-2286.1181496617564,-2264.697404030584]', '[-158.08098103415045,-86.75434269438702]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.015524148941040039, '[0.006371473312929187,0.006439419774935288]', 0, 88, 1, 800),
    @(8, 5, 1, 25, 0, 100, 0, 252, 0, '1000P periods', 0.01, '[0.0,2.4604275010558124e-11]', '[[0.0,0.0,1484.9769079581934],[0.0,0.0,0.0]]', '[-2285.005000838378,-2258.9631270010527]', '[-162.51945969096712,-82.68058198699418]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.016978979110717773, '[0.00637238166606447,0.006450803221022618]', 0, 89, 1, 902),
    @(8, 5, 1, 25, 0, 100, 0, 270, 0, '1000P periods', 0.01, '[0.0,1.237748332547009e-7]', '[[0.0,0.0,1124.8520681243645],[0.0,0.0,0.0]]', '[-2285.201589598474,-2253.9102714264836]', '[-169.53155531179323,-65.42927959465982]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 71.28549599647522, '[0.006347217617245445,0.006470293264954129]', 0, 90, 1, 4826105),
    @(8, 5, 1, 25, 0, 100, 0, 288, 0, '1000P periods', 0.01, '[-5.048162638722419e-12,4.486952003121774e-11]', '[[0.0,0.0,1484.9769079265357],[0.0,0.0,0.0]]', '[-2318.57575146697,-2210.1568775996516]', '[-205.1847662041149,8.908601259719717]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.037629127502441406, '[0.004172736119501224,0.006719477100050219]', 0, 91, 1, 1886),
    @(8, 5, 1, 25, 0, 100, 0, 306, 0, '1000P periods', 0.01, '[-1.263498438675084e-8,2.435542761871689e-9]', '[[0.0,0.0,2427.795450201531],[0.0,0.0,0.0]]', '[-4396.694968905555,3961.632725953993]', '[-178.53358453896806,79593.0331922949]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.0482330322265625, '[3.685192182672663e-6,0.016408017683696045]', 0, 92, 1, 2470),
    @(8, 5, 1, 25, 0, 100, 0, 324, 0, '1000P periods', 0.01, '[-5.453272901181741e-9,6.509055503606367e-12]', '[[0.0,0.0,3593.1832587224185],[0.0,0.0,3593.1832588825578]]', '[-4130.613873885532,-644.2667513599572]', '[-12709.623409818627,31033.456904803228]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.0536189079284668, '[2.4410402403405182e-5,0.01221844492107093]', 0, 93, 1, 2996),
    @(8, 5, 1, 25, 0, 100, 0, 342, 0, '1000P periods', 0.01, '[-2.041130984977972e-10,1.299985655318103e-10]', '[[0.0,0.0,4536.001801003193],[0.0,0.0,4536.001801115255]]', '[-2541.2525393865217,-2032.2190949835544]', '[-336.83500835198697,329.25982684066446]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', '[[0.0,0.0,0.0],[0.0,0.0,0.0]]', 0.03743696212768555, '[0.0023328370399377293,0.011461353231489195]', 0, 94, 1, 1900)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
